# Update the "想去人数" (number of people interested) figures that were
# refreshed by the generator run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9343
$ws1.Range("F4").Value = 21
$ws1.Range("F6").Value = 461

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9343
$ws4.Range("F4").Value = 21
$ws4.Range("F7").Value = 461
